$wb = $excel.ActiveWorkbook

# Sheet order: 1 Funciones_Objetivo, 2 Restricciones_del_lider,
# 3 Restricciones_del_follower, 4 Punto_modificado, 5 Vector_bf, 6 Vector_BF, 7 Vector_Alpha

# --- Sheet: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item(3)

# Row 2
$ws.Range("A2").Value = "1.5092421543248782y_1 + 0.3203062014027182y_2"
$ws.Range("B2").Value = 7.706172739445797
$ws.Range("C2").Value = "J_0_L0_v"
$ws.Range("D2").Value = 0.24011722556595838
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.8626017230215338

# Row 3
$ws.Range("A3").Value = "-4 + 0.925804466693197y_1 - 0.009471102417751735y_2"
$ws.Range("B3").Value = 0.29210593221048425
$ws.Range("C3").Value = "J_0_L0_v"
$ws.Range("D3").Value = 0.1083236165390392
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.5145103609416063

# Row 4
$ws.Range("A4").Value = "-16 - 2x + 2.240041564621893y_1 - 0.35231027860011516y_2"
$ws.Range("B4").Value = -18.061482347416565
$ws.Range("C4").Value = "J_0_LP_v"
$ws.Range("D4").Value = 0.1102758390135593
$ws.Range("E4").Value = 0.2686673865651241
$ws.Range("F4").Value = 0

# Row 5
$ws.Range("A5").Value = "-48 + 8x + 0.016033409943699617y_1 - 0.12560390005599806y_2"
$ws.Range("B5").Value = -1.1877764405453983
$ws.Range("C5").Value = "J_Ne_L0_v"
$ws.Range("D5").Value = 0.3168885247170169
$ws.Range("E5").Value = 0.8383936992491441
$ws.Range("F5").Value = 0

# Row 6
$ws.Range("A6").Value = "12 - 2x - 0.032066819887399234y_1 + 0.2512078001119961y_2"
$ws.Range("B6").Value = 0.6305578800946663
$ws.Range("C6").Value = "J_Ne_L0_v"
$ws.Range("D6").Value = 0.4167665579899481
$ws.Range("E6").Value = 0.29544899638169286
$ws.Range("F6").Value = 0

# --- Sheet: Punto_modificado ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 5.875357499928848
$ws.Range("B2").Value = 4.657691821664619
$ws.Range("C2").Value = 2.112315956957238

# --- Sheet: Vector_bf (lowercase bf) ---
# NOTE: worksheet names "Vector_bf" and "Vector_BF" differ only by case and
# Excel sheet name lookups are case-insensitive, so using Worksheets.Item("Vector_BF")
# would ambiguously resolve to the first ("Vector_bf"). Use positional index instead:
# sheet order: 1 Funciones_Objetivo, 2 Restricciones_del_lider,
# 3 Restricciones_del_follower, 4 Punto_modificado, 5 Vector_bf, 6 Vector_BF, 7 Vector_Alpha
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 0.2985795847191042
$ws.Range("A3").Value = -0.1019263563856154

# --- Sheet: Vector_BF (uppercase BF) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = -4.578916828099519
$ws.Range("A3").Value = 2.394205686914393
$ws.Range("A4").Value = -1.8742592922067505

# --- Sheet: Vector_Alpha ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 0.5564707754203111
$ws.Range("A3").Value = 0.07103381391839458
